$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Actualizacion_de_Datos" rows (24-26), copying the
# formatting (style s="5") of B23 down into each new row.
$ws.Range("B23").Copy($ws.Range("B24:B26"))

$ws.Range("A24").Value = "Actualizacion_de_Datos_MIX"
$ws.Range("B24").Value = 1162816939

$ws.Range("A25").Value = "Actualizacion_de_Datos_POS"
$ws.Range("B25").Value = 1145642605

$ws.Range("A26").Value = "Actualizacion_de_Datos_PRE"
$ws.Range("B26").Value = 1162676705

# Update the sheet view: move the selection to C22 (also clears the old
# topLeftCell="A3" scroll position, matching the saved view).
$ws.Activate()
$ws.Range("C22").Select()
